# Add a "time_taken" column (F) with the header styled like the other
# header cells (A1:E1) and per-row timestamp values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (E1) onto F1 so the new
# header reuses the same style (bold, centered, bordered) instead of
# creating a brand-new style entry.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "time_taken"

$ws.Range("F2").Value = "2021-10-05 13:40:19.721132"
$ws.Range("F3").Value = "2021-10-05 13:40:19.721143"
$ws.Range("F4").Value = "2021-10-05 13:40:19.721147"
